$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 81 (shifts existing rows 81-91 down to 85-95)
$ws.Rows.Item(81).Resize(4).Insert()

# Fill the 4 new rows (81-84) with a new week of "Chirimoya" price data
# for "Provincia de Limarí" (fecha serial 44461).

# Row 81: Especial
$ws.Range("A81").Value = 6
$ws.Range("B81").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value = 44461
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100107
$ws.Range("H81").Value = "Otros"
$ws.Range("I81").Value = 100107002
$ws.Range("J81").Value = "Chirimoya"
$ws.Range("K81").Value = "Cultivar IV Región"
$ws.Range("L81").Value = "Especial"
$ws.Range("M81").Value = 200
$ws.Range("N81").Value = 3000
$ws.Range("O81").Value = 3000
$ws.Range("P81").Value = 3000
$ws.Range("Q81").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R81").Value = "Provincia de Limarí"
$ws.Range("S81").Value = 3000
$ws.Range("T81").Value = 1

# Row 82: Extra (doble especial)
$ws.Range("A82").Value = 6
$ws.Range("B82").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C82").Value = "Metropolitana"
$ws.Range("D82").Value = 44461
$ws.Range("E82").Value = 13
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100107
$ws.Range("H82").Value = "Otros"
$ws.Range("I82").Value = 100107002
$ws.Range("J82").Value = "Chirimoya"
$ws.Range("K82").Value = "Cultivar IV Región"
$ws.Range("L82").Value = "Extra (doble especial)"
$ws.Range("M82").Value = 150
$ws.Range("N82").Value = 3200
$ws.Range("O82").Value = 3200
$ws.Range("P82").Value = 3200
$ws.Range("Q82").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R82").Value = "Provincia de Limarí"
$ws.Range("S82").Value = 3200
$ws.Range("T82").Value = 1

# Row 83: Primera
$ws.Range("A83").Value = 6
$ws.Range("B83").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C83").Value = "Metropolitana"
$ws.Range("D83").Value = 44461
$ws.Range("E83").Value = 13
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100107
$ws.Range("H83").Value = "Otros"
$ws.Range("I83").Value = 100107002
$ws.Range("J83").Value = "Chirimoya"
$ws.Range("K83").Value = "Cultivar IV Región"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 270
$ws.Range("N83").Value = 2500
$ws.Range("O83").Value = 2600
$ws.Range("P83").Value = 2550
$ws.Range("Q83").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R83").Value = "Provincia de Limarí"
$ws.Range("S83").Value = 2550
$ws.Range("T83").Value = 1

# Row 84: Segunda
$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44461
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100107
$ws.Range("H84").Value = "Otros"
$ws.Range("I84").Value = 100107002
$ws.Range("J84").Value = "Chirimoya"
$ws.Range("K84").Value = "Cultivar IV Región"
$ws.Range("L84").Value = "Segunda"
$ws.Range("M84").Value = 270
$ws.Range("N84").Value = 2000
$ws.Range("O84").Value = 2000
$ws.Range("P84").Value = 2000
$ws.Range("Q84").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R84").Value = "Provincia de Limarí"
$ws.Range("S84").Value = 2000
$ws.Range("T84").Value = 1
